# Add new Korean/English phrase rows (32-51) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('자료 준비되었습니다.', 'The materials are ready.')
    ,@('언제 시간 되시나요?', 'When are you available?')
    ,@('점심 식사 하셨습니까?', 'Have you had lunch?')
    ,@('커피 한 잔 하시겠습니까?', 'Would you like a cup of coffee?')
    ,@('잠시 자리 비우셨습니다.', 'He/She is away from the desk for a moment.')
    ,@('메모 남겨드릴까요?', 'Shall I leave a message?')
    ,@('팩스로 보내주세요.', 'Please send it by fax.')
    ,@('프레젠테이션 준비 완료했습니다.', 'Presentation preparation is complete.')
    ,@('예산안 검토 부탁드립니다.', 'Please review the budget proposal.')
    ,@('다음 주 일정 잡겠습니다.', 'I will schedule it for next week.')
    ,@('출근했습니다.', 'I have arrived at work.')
    ,@('병가 내겠습니다.', 'I will take a sick leave.')
    ,@('연차 쓰겠습니다.', 'I will use my annual leave.')
    ,@('급한 일입니까?', 'Is it urgent?')
    ,@('우선순위가 어떻게 되나요?', 'What is the priority?')
    ,@('피드백 부탁드립니다.', 'Please give me feedback.')
    ,@('수정해서 다시 보내겠습니다.', 'I will revise and resend it.')
    ,@('성공적인 프로젝트였습니다.', 'It was a successful project.')
    ,@('함께 일해서 즐거웠습니다.', 'It was a pleasure working with you.')
    ,@('앞으로도 잘 부탁드립니다.', 'I look forward to working with you.')
)

$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

